$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new primary-key column before column A; existing data shifts right.
$ws.Columns("A:A").Insert()

# Fill in the new primary key column (row numbers 1..6), centered.
$ws.Range("A1").Value = 1
$ws.Range("A2").Value = 2
$ws.Range("A3").Value = 3
$ws.Range("A4").Value = 4
$ws.Range("A5").Value = 5
$ws.Range("A6").Value = 6

$idRange = $ws.Range("A1:A6")
$idRange.HorizontalAlignment = -4108
$idRange.VerticalAlignment = -4108
$ws.Columns("A:A").ColumnWidth = 4.5703125

# Re-apply a uniform thin border around every data cell (replacing the old
# medium partial border) across the whole table, including the new column.
$allRange = $ws.Range("A1:G6")
$allRange.Borders.LineStyle = -4142
$allRange.Borders.LineStyle = 1
$allRange.Borders.Weight = 2

# Restore the active selection cell noted in the authored edit.
$null = $ws.Range("J2").Select()
